$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at the top of the data block (row 661), pushing the
# existing rows 661:689 down to 665:693.
$ws.Rows("661:664").Insert()

# Row 661 - new weekly entry (Provincia de Quillota, Primera)
$ws.Range("A661").Value = 11
$ws.Range("B661").Value = "Vega Monumental Concepción"
$ws.Range("C661").Value = "Bíobío"
$ws.Range("D661").Value = 45008
$ws.Range("E661").Value = 8
$ws.Range("F661").Value = 100112020
$ws.Range("G661").Value = "Tomate"
$ws.Range("H661").Value = "Larga vida"
$ws.Range("I661").Value = "Primera"
$ws.Range("J661").Value = 300
$ws.Range("K661").Value = 12000
$ws.Range("L661").Value = 12000
$ws.Range("M661").Value = 12000
$ws.Range("N661").Value = "$/bandeja 18 kilos"
$ws.Range("O661").Value = "Provincia de Quillota"
$ws.Range("P661").Value = 667
$ws.Range("Q661").Value = 18
$ws.Range("R661").Value = "Hortaliza"

# Row 662 - new weekly entry (Provincia de Quillota, Segunda)
$ws.Range("A662").Value = 11
$ws.Range("B662").Value = "Vega Monumental Concepción"
$ws.Range("C662").Value = "Bíobío"
$ws.Range("D662").Value = 45008
$ws.Range("E662").Value = 8
$ws.Range("F662").Value = 100112020
$ws.Range("G662").Value = "Tomate"
$ws.Range("H662").Value = "Larga vida"
$ws.Range("I662").Value = "Segunda"
$ws.Range("J662").Value = 300
$ws.Range("K662").Value = 10000
$ws.Range("L662").Value = 10000
$ws.Range("M662").Value = 10000
$ws.Range("N662").Value = "$/bandeja 18 kilos"
$ws.Range("O662").Value = "Provincia de Quillota"
$ws.Range("P662").Value = 556
$ws.Range("Q662").Value = 18
$ws.Range("R662").Value = "Hortaliza"

# Row 663 - new weekly entry (Región Metropolitana, Semiduro, Primera)
$ws.Range("A663").Value = 11
$ws.Range("B663").Value = "Vega Monumental Concepción"
$ws.Range("C663").Value = "Bíobío"
$ws.Range("D663").Value = 45008
$ws.Range("E663").Value = 8
$ws.Range("F663").Value = 100112020
$ws.Range("G663").Value = "Tomate"
$ws.Range("H663").Value = "Semiduro"
$ws.Range("I663").Value = "Primera"
$ws.Range("J663").Value = 500
$ws.Range("K663").Value = 7000
$ws.Range("L663").Value = 7500
$ws.Range("M663").Value = 7200
$ws.Range("N663").Value = "$/bandeja 18 kilos"
$ws.Range("O663").Value = "Región Metropolitana"
$ws.Range("P663").Value = 400
$ws.Range("Q663").Value = 18
$ws.Range("R663").Value = "Hortaliza"

# Row 664 - new weekly entry (Región Metropolitana, Semiduro, Segunda)
$ws.Range("A664").Value = 11
$ws.Range("B664").Value = "Vega Monumental Concepción"
$ws.Range("C664").Value = "Bíobío"
$ws.Range("D664").Value = 45008
$ws.Range("E664").Value = 8
$ws.Range("F664").Value = 100112020
$ws.Range("G664").Value = "Tomate"
$ws.Range("H664").Value = "Semiduro"
$ws.Range("I664").Value = "Segunda"
$ws.Range("J664").Value = 300
$ws.Range("K664").Value = 6000
$ws.Range("L664").Value = 6000
$ws.Range("M664").Value = 6000
$ws.Range("N664").Value = "$/bandeja 18 kilos"
$ws.Range("O664").Value = "Región Metropolitana"
$ws.Range("P664").Value = 333
$ws.Range("Q664").Value = 18
$ws.Range("R664").Value = "Hortaliza"

# Match the date-format style used by the rest of column D for the new rows
$ws.Range("D661:D664").NumberFormat = $ws.Range("D665").NumberFormat
